$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.054.39'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.38%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.828.16'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.31%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9987'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '240.59'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.28%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6185'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -7.00%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.000'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '44.45'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +6.30%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07342'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.10%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.2915'
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '22.68'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.19%  '
$ws.Range("E12").Value = '  -0.62%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.825.56'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.41%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.962'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.47%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6618'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.08%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '81.78'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.45%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000009014'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +7.30%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.027'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.20%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '29.038.85'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.21%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '2.073.78'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.35%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '225.56'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.78%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '12.36'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.83%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.001'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.04%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.120'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.79%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.0000'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.01%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '160.02'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.25%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.425'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.38%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.1354'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -4.23%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '17.78'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.01%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.492'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.39%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.041'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.12%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.051'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.58%  '
$ws.Range("E33").Value = '  +0.35%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05264'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.25%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.839'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.22%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.149'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.04%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.7313'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -3.63%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.648'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.89%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.300.84'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.10%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.749'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.52%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.01781'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.98%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.296'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +5.37%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.9030'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.89%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.9992'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.63%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '101.84'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.05%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.972.06'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.79%  '
$ws.Range("B47").Value = 'BabyDogeCoin'
$ws.Range("C47").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.00000000122'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.76%  '
$ws.Range("B48").Value = 'Mantle'
$ws.Range("C48").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.5114'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.97%  '
$ws.Range("B49").Value = 'Aave'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '63.88'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.70%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.716'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.20%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.3963'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.81%  '
